# Correct mistake in trad gravity - ppml
# Updates the "(4) PPML" column of the "Traditional Gravity Estimates"
# table (Table 1) with corrected coefficient / SE / N / logLik / AIC values.
#
# Note: cells are addressed directly via Table.Cell(row, col) (1-based, as
# in the Word object model) and the replacement text is written straight
# onto the cell's Range so the edit stays scoped to that single cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    # Cell.Range.Text includes the trailing cell-mark/paragraph-mark
    # characters (CR + cell marker), so trim those before comparing.
    $current = $cell.Range.Text.TrimEnd([char]0x0D, [char]0x07)
    if ($current -ne $oldText) {
        throw "Unexpected text '$current' in table cell ($row, $col); expected '$oldText'"
    }
    $cell.Range.Text = $newText
}

# Contiguity coefficient: 0.438 -> 0.437
Set-CellValue $t 7 5 "0.438" "0.437"

# Contiguity standard error: (0.085) -> (0.084)
Set-CellValue $t 8 5 "(0.085)" "(0.084)"

# Common language coefficient: 0.246 -> 0.247
Set-CellValue $t 9 5 "0.246" "0.247"

# Colony coefficient: -0.223 -> -0.222
Set-CellValue $t 11 5 "-0.223" "-0.222"

# N: 25689 -> 28152
Set-CellValue $t 21 5 "25689" "28152"

# logLik: -2182850 -> -2194537
Set-CellValue $t 23 5 "-2182850" "-2194537"

# AIC: 4367361 -> 4390736
Set-CellValue $t 24 5 "4367361" "4390736"
